# Updated cryptos list (Price in column D, Volume(1h) in column E)
# For cells whose new value is a plain decimal number (e.g. "6.01"), the
# NumberFormat is forced to Text ("@") first so Excel keeps storing the
# cell as a text string (matching the source data, e.g. prices like
# "37.333.63" which aren't valid numbers) instead of auto-converting it
# to a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "37.333.63"
$ws.Cells.Item(2, 5).Value = "  -1.30%  "
$ws.Cells.Item(3, 4).Value = "2.051.23"
$ws.Cells.Item(3, 5).Value = "  -1.26%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "230.56"
$ws.Cells.Item(5, 5).Value = "  -1.31%  "
$ws.Cells.Item(6, 5).Value = "  -0.76%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "57.07"
$ws.Cells.Item(8, 5).Value = "  -3.75%  "
$ws.Cells.Item(9, 5).Value = "  -2.31%  "
$ws.Cells.Item(10, 5).Value = "  -2.67%  "
$ws.Cells.Item(11, 5).Value = "  +1.29%  "
$ws.Cells.Item(12, 5).Value = "  -0.58%  "
$ws.Cells.Item(13, 4).Value = "2.354.37"
$ws.Cells.Item(13, 5).Value = "  -1.27%  "
$ws.Cells.Item(14, 5).Value = "  -3.06%  "
$ws.Cells.Item(15, 5).Value = "  -2.42%  "
$ws.Cells.Item(16, 5).Value = "  -2.11%  "
$ws.Cells.Item(17, 4).Value = "2.039.52"
$ws.Cells.Item(17, 5).Value = "  -2.05%  "
$ws.Cells.Item(18, 4).Value = "37.297.77"
$ws.Cells.Item(18, 5).Value = "  -1.11%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.01"
$ws.Cells.Item(19, 5).Value = "  -2.33%  "
$ws.Cells.Item(20, 5).Value = "  -2.51%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0823"
$ws.Cells.Item(21, 5).Value = "  -3.78%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "226.68"
$ws.Cells.Item(22, 5).Value = "  -0.72%  "
$ws.Cells.Item(23, 5).Value = "  +0.07%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.38"
$ws.Cells.Item(24, 5).Value = "  +0.65%  "
$ws.Cells.Item(25, 5).Value = "  -3.77%  "
$ws.Cells.Item(26, 5).Value = "  +3.04%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "168.80"
$ws.Cells.Item(27, 5).Value = "  -1.03%  "
$ws.Cells.Item(28, 5).Value = "  -3.33%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "19.14"
$ws.Cells.Item(29, 5).Value = "  -1.84%  "
$ws.Cells.Item(30, 5).Value = "  -5.18%  "
$ws.Cells.Item(31, 5).Value = "  +0.11%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.52"
$ws.Cells.Item(32, 5).Value = "  -4.20%  "
$ws.Cells.Item(33, 5).Value = "  -1.45%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.55"
$ws.Cells.Item(35, 5).Value = "  -0.07%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.82"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "
$ws.Cells.Item(37, 5).Value = "  -3.83%  "
$ws.Cells.Item(38, 5).Value = "  -0.09%  "
$ws.Cells.Item(39, 5).Value = "  -2.61%  "
$ws.Cells.Item(40, 5).Value = "  +4.54%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "98.07"
$ws.Cells.Item(41, 5).Value = "  -1.28%  "
$ws.Cells.Item(42, 4).Value = "1.485.34"
$ws.Cells.Item(43, 5).Value = "  +0.74%  "
$ws.Cells.Item(44, 5).Value = "  -3.49%  "
$ws.Cells.Item(45, 5).Value = "  +2.90%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "16.60"
$ws.Cells.Item(46, 5).Value = "  -0.77%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.02"
$ws.Cells.Item(47, 5).Value = "  -2.93%  "
$ws.Cells.Item(48, 5).Value = "  -3.31%  "
$ws.Cells.Item(49, 5).Value = "  -2.35%  "
$ws.Cells.Item(50, 5).Value = "  -2.85%  "
$ws.Cells.Item(51, 4).Value = "2.239.84"
$ws.Cells.Item(51, 5).Value = "  -1.31%  "
